$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

function Set-TextForceNumericLooking($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-Text "D2" "58.084.17"
Set-Text "E2" "  +1.49%  "

# Row 3 - Ethereum
Set-Text "D3" "3.137.78"
Set-Text "E3" "  +1.42%  "

# Row 4 - TetherUSD
Set-Text "E4" "  +0.00%  "

# Row 5 - BNB
Set-TextForceNumericLooking "D5" "535.27"
Set-Text "E5" "  +2.28%  "

# Row 6 - Solana
Set-TextForceNumericLooking "D6" "139.17"
Set-Text "E6" "  +1.97%  "

# Row 7 - USDC
Set-Text "E7" "  +0.09%  "

# Row 8 - XRP
Set-Text "E8" "  +11.37%  "

# Row 9 - Toncoin
Set-TextForceNumericLooking "D9" "7.33"
Set-Text "E9" "  -0.30%  "

# Row 10 - Dogecoin
Set-Text "E10" "  +2.42%  "

# Row 11 - Cardano
Set-TextForceNumericLooking "D11" "0.421"
Set-Text "E11" "  +5.82%  "

# Row 12 - TRON
Set-Text "E12" "  +3.59%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-Text "D13" "3.679.41"
Set-Text "E13" "  +1.45%  "

# Row 14 - Avalanche
Set-TextForceNumericLooking "D14" "25.73"
Set-Text "E14" "  +1.41%  "

# Row 15 - ShibaInu
Set-TextForceNumericLooking "D15" "0.0000170"
Set-Text "E15" "  +4.91%  "

# Row 16 - WrappedBTC
Set-Text "D16" "58.147.01"
Set-Text "E16" "  +1.44%  "

# Row 17 - Polkadot
Set-Text "E17" "  +6.02%  "

# Row 18 - WrappedEther
Set-Text "D18" "3.139.96"
Set-Text "E18" "  +1.39%  "

# Row 19 - Chainlink
Set-TextForceNumericLooking "D19" "12.99"
Set-Text "E19" "  +4.16%  "

# Row 20 - Uniswap
Set-TextForceNumericLooking "D20" "8.21"
Set-Text "E20" "  +4.49%  "

# Row 21 - BitcoinCash
Set-TextForceNumericLooking "D21" "377.01"
Set-Text "E21" "  +8.02%  "

# Row 22 - Dai
Set-Text "E22" "  +0.08%  "

# Row 23 - LEO
Set-TextForceNumericLooking "D23" "5.72"
Set-Text "E23" "  -1.00%  "

# Row 24 - Litecoin
Set-TextForceNumericLooking "D24" "70.16"
Set-Text "E24" "  +2.25%  "

# Row 25 - Polygon
Set-TextForceNumericLooking "D25" "0.516"
Set-Text "E25" "  +3.30%  "

# Row 26 - Kaspa
Set-TextForceNumericLooking "D26" "0.167"
Set-Text "E26" "  +0.38%  "

# Row 27 - Binance-PegBSC-USD
Set-TextForceNumericLooking "D27" "0.996"
Set-Text "E27" "  -0.18%  "

# Row 28 - InternetComputer(DFINITY) [swapped into row28]
Set-Text "B28" "InternetComputer(DFINITY)"
Set-Text "C28" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextForceNumericLooking "D28" "7.99"
Set-Text "E28" "  +10.34%  "

# Row 29 - PEPE [swapped into row29]
Set-Text "B29" "PEPE"
Set-Text "C29" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-Text "D29" "0.0₃0885"
Set-Text "E29" "  +1.64%  "

# Row 30 - RenderToken
Set-TextForceNumericLooking "D30" "6.18"
Set-Text "E30" "  +5.13%  "

# Row 31 - PancakeSwap
Set-TextForceNumericLooking "D31" "1.89"
Set-Text "E31" "  +1.03%  "

# Row 32 - EthereumClassic
Set-TextForceNumericLooking "D32" "21.77"
Set-Text "E32" "  +4.16%  "

# Row 33 - NEARProtocol
Set-TextForceNumericLooking "D33" "5.19"
Set-Text "E33" "  +5.83%  "

# Row 34 - Fetch.AI
Set-TextForceNumericLooking "D34" "1.18"
Set-Text "E34" "  +3.02%  "

# Row 35 - Monero
Set-TextForceNumericLooking "D35" "161.38"
Set-Text "E35" "  +1.45%  "

# Row 36 - Aptos
Set-TextForceNumericLooking "D36" "6.29"
Set-Text "E36" "  +4.21%  "

# Row 37 - ImmutableX
Set-Text "E37" "  +8.62%  "

# Row 38 - EnergySwap
Set-TextForceNumericLooking "D38" "25.59"
Set-Text "E38" "  -0.23%  "

# Row 39 - Stacks
Set-Text "E39" "  +5.09%  "

# Row 40 - Maker
Set-Text "D40" "2.632.32"
Set-Text "E40" "  +9.53%  "

# Row 41 - Filecoin [swapped into row41]
Set-Text "B41" "Filecoin"
Set-Text "C41" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextForceNumericLooking "D41" "4.26"
Set-Text "E41" "  +5.42%  "

# Row 42 - Hedera [swapped into row42]
Set-Text "B42" "Hedera"
Set-Text "C42" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextForceNumericLooking "D42" "0.0675"
Set-Text "E42" "  +2.30%  "

# Row 43 - OKB
Set-TextForceNumericLooking "D43" "38.98"
Set-Text "E43" "  +6.14%  "

# Row 44 - Mantle
Set-TextForceNumericLooking "D44" "0.701"
Set-Text "E44" "  +0.70%  "

# Row 45 - VeChain
Set-TextForceNumericLooking "D45" "0.0272"
Set-Text "E45" "  +3.59%  "

# Row 46 - FirstDigitalUSD
Set-Text "E46" "  -0.04%  "

# Row 47 - Cosmos
Set-TextForceNumericLooking "D47" "6.24"
Set-Text "E47" "  +4.49%  "

# Row 48 - ONDO
Set-TextForceNumericLooking "D48" "0.979"
Set-Text "E48" "  +2.35%  "

# Row 49 - Stellar
Set-TextForceNumericLooking "D49" "0.100"
Set-Text "E49" "  +10.50%  "

# Row 50 - InjectiveProtocol
Set-TextForceNumericLooking "D50" "20.35"
Set-Text "E50" "  +2.82%  "

# Row 51 - SuiNetwork
Set-Text "E51" "  -1.63%  "
